# Update "想去人数" (F column) counts on the 展览 (sheet1), 演出 (sheet2)
# and 全部类型 (sheet4) worksheets to the values captured at the latest
# gh-pages data refresh (commit 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 239
$ws1.Range("F4").Value = 4871
$ws1.Range("F5").Value = 213
$ws1.Range("F6").Value = 170
$ws1.Range("F7").Value = 129
$ws1.Range("F8").Value = 120
$ws1.Range("F9").Value = 101
$ws1.Range("F10").Value = 779
$ws1.Range("F11").Value = 242
$ws1.Range("F12").Value = 1238
$ws1.Range("F13").Value = 131
$ws1.Range("F14").Value = 258
$ws1.Range("F15").Value = 208
$ws1.Range("F16").Value = 91
$ws1.Range("F18").Value = 159
$ws1.Range("F19").Value = 118
$ws1.Range("F20").Value = 4185
$ws1.Range("F21").Value = 6488
$ws1.Range("F25").Value = 546
$ws1.Range("F26").Value = 49
$ws1.Range("F27").Value = 4019
$ws1.Range("F28").Value = 416
$ws1.Range("F31").Value = 2621
$ws1.Range("F33").Value = 537
$ws1.Range("F34").Value = 154
$ws1.Range("F35").Value = 312
$ws1.Range("F36").Value = 330
$ws1.Range("F37").Value = 385
$ws1.Range("F38").Value = 196
$ws1.Range("F39").Value = 15
$ws1.Range("F40").Value = 1583
$ws1.Range("F41").Value = 989
$ws1.Range("F42").Value = 53
$ws1.Range("F43").Value = 90
$ws1.Range("F44").Value = 62
$ws1.Range("F45").Value = 505
$ws1.Range("F46").Value = 487
$ws1.Range("F48").Value = 81
$ws1.Range("F49").Value = 600

# --- Sheet "演出" ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 113

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 239
$ws4.Range("F4").Value = 4871
$ws4.Range("F5").Value = 213
$ws4.Range("F6").Value = 170
$ws4.Range("F7").Value = 129
$ws4.Range("F8").Value = 113
$ws4.Range("F9").Value = 120
$ws4.Range("F10").Value = 101
$ws4.Range("F11").Value = 779
$ws4.Range("F12").Value = 242
$ws4.Range("F13").Value = 1238
$ws4.Range("F14").Value = 131
$ws4.Range("F15").Value = 208
$ws4.Range("F16").Value = 91
$ws4.Range("F18").Value = 159
$ws4.Range("F19").Value = 118
$ws4.Range("F20").Value = 4185
$ws4.Range("F21").Value = 6488
$ws4.Range("F25").Value = 546
$ws4.Range("F26").Value = 49
$ws4.Range("F27").Value = 4019
$ws4.Range("F28").Value = 416
$ws4.Range("F31").Value = 2621
$ws4.Range("F33").Value = 537
$ws4.Range("F34").Value = 154
$ws4.Range("F35").Value = 312
$ws4.Range("F36").Value = 330
$ws4.Range("F37").Value = 385
$ws4.Range("F38").Value = 196
$ws4.Range("F39").Value = 15
$ws4.Range("F40").Value = 1583
$ws4.Range("F41").Value = 989
$ws4.Range("F42").Value = 53
$ws4.Range("F43").Value = 90
$ws4.Range("F44").Value = 62
$ws4.Range("F45").Value = 505
$ws4.Range("F46").Value = 487
$ws4.Range("F48").Value = 81
$ws4.Range("F49").Value = 600
